# The workbook tracks daily/weekly price observations for
# "Pepino ensalada" at Terminal La Palmera de La Serena. This edit adds
# two brand-new observations (dated 45131) at the top of the data block
# (rows 825-826), pushing the previously-existing rows 825-932 down by
# two positions (to rows 827-934). The sheet's used-range grows from
# A1:R932 to A1:R934 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 825 - this
# shifts all data from row 825 downward by two rows (825->827 ... 932->934)
# while Excel automatically carries the date-format style (column D)
# down with the shifted cells.
$ws.Rows("825:826").Insert()

# --- New row 825: "Primera" quality entry for date 45131 ---
$ws.Cells(825,1).Value2 = 8
$ws.Cells(825,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells(825,3).Value2 = "Coquimbo"
$ws.Cells(825,4).Value2 = 45131
$ws.Cells(825,5).Value2 = 4
$ws.Cells(825,6).Value2 = 100112043
$ws.Cells(825,7).Value2 = "Pepino ensalada"
$ws.Cells(825,8).Value2 = "Sin especificar"
$ws.Cells(825,9).Value2 = "Primera"
$ws.Cells(825,10).Value2 = 520
$ws.Cells(825,11).Value2 = 9000
$ws.Cells(825,12).Value2 = 10000
$ws.Cells(825,13).Value2 = 9500
$ws.Cells(825,14).Value2 = "$/caja 60 unidades"
$ws.Cells(825,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells(825,16).Value2 = 158
$ws.Cells(825,17).Value2 = 60
$ws.Cells(825,18).Value2 = "Hortaliza"

# --- New row 826: "Segunda" quality entry for date 45131 ---
$ws.Cells(826,1).Value2 = 8
$ws.Cells(826,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells(826,3).Value2 = "Coquimbo"
$ws.Cells(826,4).Value2 = 45131
$ws.Cells(826,5).Value2 = 4
$ws.Cells(826,6).Value2 = 100112043
$ws.Cells(826,7).Value2 = "Pepino ensalada"
$ws.Cells(826,8).Value2 = "Sin especificar"
$ws.Cells(826,9).Value2 = "Segunda"
$ws.Cells(826,10).Value2 = 320
$ws.Cells(826,11).Value2 = 6000
$ws.Cells(826,12).Value2 = 7000
$ws.Cells(826,13).Value2 = 6500
$ws.Cells(826,14).Value2 = "$/caja 80 unidades"
$ws.Cells(826,15).Value2 = "Región de Arica y Parinacota"
$ws.Cells(826,16).Value2 = 81
$ws.Cells(826,17).Value2 = 80
$ws.Cells(826,18).Value2 = "Hortaliza"
